# Large-scale mode: scale the intervention results by 5373x
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scale = 5373

for ($r = 2; $r -le 22; $r++) {
    foreach ($col in @("B", "C")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne 0) {
            $cell.Value2 = $val * $scale
        }
    }
}
